# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown -
# correct the "Date" column (BF) from the malformed "6-4-2012-13"
# text to the proper ISO date string "2013-06-04" for every data row.
#
# The value must remain literal text (not get auto-converted into a
# date serial by Excel's smart input parsing), so it is entered with a
# leading apostrophe (text qualifier) via Formula, exactly as typing
# '2013-06-04 into the cell would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "6-4-2012-13") {
        $cell.Formula = "'2013-06-04"
    }
}
